$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global_Variables")

# F2 becomes a numeric 1 instead of the text "1,2,3"
$ws.Range("F2").Value = 1

# F3 / F4 gain a numeric 1 as well
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 1

# New header row for the two extra paystub option columns
$ws.Range("H1").Value = "paystub_B_options"
$ws.Range("I1").Value = "paystub_C_options"

# paystub_A_options (column G) values per row
$ws.Range("G2").Value = '{"Rate" : 20 , "4_Digit_Account_Number" : 8698, "Numbe of Paystubs" : 5, "Period" : "Apr 01 2022"}'
$ws.Range("G3").Value = '{"Rate" : 25 , "4_Digit_Account_Number" : 8698, "Numbe of Paystubs" : 5, "Period" : "Apr 01 2022"}'
$ws.Range("G4").Value = '{"Rate" : 30 , "4_Digit_Account_Number" : 8698, "Numbe of Paystubs" : 3, "Period" : "Apr 01 2022"}'

# paystub_B_options (column H) values per row
$ws.Range("H2").Value = '{"Rate" : 20 , "occupation" : "Student", "Numbe of Paystubs" : 5, "Period" : "Apr 01 2022"}'
$ws.Range("H3").Value = '{"Rate" : 20 , "occupation" : "Student", "Numbe of Paystubs" : 5, "Period" : "Apr 01 2022"}'
$ws.Range("H4").Value = '{"Rate" : 20 , "occupation" : "Student", "Numbe of Paystubs" : 5, "Period" : "Apr 01 2022"}'

# paystub_C_options (column I) values per row
$ws.Range("I2").Value = '{"Rate" : 20 , "Numbe of Paystubs" : 5, "Period" : "Apr 01 2022"}'
$ws.Range("I3").Value = '{"Rate" : 20 , "Numbe of Paystubs" : 5, "Period" : "Apr 01 2022"}'
$ws.Range("I4").Value = '{"Rate" : 20 , "Numbe of Paystubs" : 5, "Period" : "Apr 01 2022"}'

# Match the saved selection state from the authored workbook
$ws.Range("I4").Select()
